# Add sort algorithm name strings to column A (rows 2-8)
# and update the benchmark numbers in columns B:K (rows 2-8)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @("Bubble", "Insertion", "Selection", "Merge", "Quick", "Shell", "Arrays.sort()")

$values = @(
    @(4789, 2478, 8317, 28434, 91841, 317414, 1172677, 4879638, 24750065, 120976040),
    @(3816, 1505, 4652, 17212, 64304, 245325, 958934, 3793146, 15048841, 60182701),
    @(3542, 1900, 5550, 15859, 54679, 161239, 561097, 2091575, 7909743, 31530299),
    @(9731, 2615, 5215, 11388, 24891, 53523, 100766, 245173, 452302, 950160),
    @(3466, 2174, 4576, 59575, 11130, 23279, 49828, 114102, 303243, 868628),
    @(2676, 1216, 3025, 7009, 16604, 38454, 89818, 206368, 474350, 1069280),
    @(9564, 6918, 2265, 3999, 14369, 20132, 38758, 69884, 131938, 256166)
)

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]

    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
